# Auto-generated edit script for cryptos.xlsx update
# Applies cell-level text updates to columns B, C, D, E per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be treated/stored as text (not auto-converted to a number/date)
    $range.NumberFormat = "@"
    $range.Value = $value
    # Reset the style back to the workbook default so we do not leave a stray
    # number-format style on the cell (matches original unstyled cells).
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.185.74"
Set-TextValue $ws.Range("E2") "  +8.94%  "
Set-TextValue $ws.Range("D3") "3.446.79"
Set-TextValue $ws.Range("E3") "  +5.59%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "413.77"
Set-TextValue $ws.Range("E5") "  +4.17%  "
Set-TextValue $ws.Range("D6") "122.79"
Set-TextValue $ws.Range("E6") "  +13.06%  "
Set-TextValue $ws.Range("D7") "3.442.09"
Set-TextValue $ws.Range("E7") "  +5.60%  "
Set-TextValue $ws.Range("D8") "0.592"
Set-TextValue $ws.Range("E8") "  +1.91%  "
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  +0.02%  "
Set-TextValue $ws.Range("D10") "0.651"
Set-TextValue $ws.Range("E10") "  +4.49%  "
Set-TextValue $ws.Range("E11") "  +30.62%  "
Set-TextValue $ws.Range("D12") "41.27"
Set-TextValue $ws.Range("E12") "  +4.54%  "
Set-TextValue $ws.Range("E13") "  -0.37%  "
Set-TextValue $ws.Range("D14") "3.972.92"
Set-TextValue $ws.Range("E14") "  +5.17%  "
Set-TextValue $ws.Range("D15") "8.49"
Set-TextValue $ws.Range("E15") "  +2.54%  "
Set-TextValue $ws.Range("D16") "19.66"
Set-TextValue $ws.Range("E16") "  +3.31%  "
Set-TextValue $ws.Range("D17") "3.431.86"
Set-TextValue $ws.Range("E17") "  +5.19%  "
Set-TextValue $ws.Range("D18") "61.945.08"
Set-TextValue $ws.Range("E18") "  +8.80%  "
Set-TextValue $ws.Range("E19") "  -0.83%  "
Set-TextValue $ws.Range("D20") "10.77"
Set-TextValue $ws.Range("E20") "  -2.44%  "
Set-TextValue $ws.Range("E21") "  +23.05%  "
Set-TextValue $ws.Range("E22") "  -1.10%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D23") "81.32"
Set-TextValue $ws.Range("E23") "  +9.32%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D24") "312.75"
Set-TextValue $ws.Range("E24") "  +6.71%  "
Set-TextValue $ws.Range("D25") "12.97"
Set-TextValue $ws.Range("E25") "  +0.23%  "
Set-TextValue $ws.Range("D26") "3.15"
Set-TextValue $ws.Range("E26") "  -0.84%  "
Set-TextValue $ws.Range("D27") "31.14"
Set-TextValue $ws.Range("E27") "  +10.50%  "
Set-TextValue $ws.Range("D28") "7.89"
Set-TextValue $ws.Range("E28") "  +6.20%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D29") "7.72"
Set-TextValue $ws.Range("E29") "  -2.74%  "
Set-TextValue $ws.Range("E30") "  +2.89%  "
$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D31") "4.29"
Set-TextValue $ws.Range("E31") "  -2.03%  "
Set-TextValue $ws.Range("E32") "  +4.06%  "
Set-TextValue $ws.Range("E33") "  +20.28%  "
Set-TextValue $ws.Range("D34") "41.96"
Set-TextValue $ws.Range("E34") "  +5.07%  "
Set-TextValue $ws.Range("E35") "  +1.56%  "
Set-TextValue $ws.Range("E36") "  +0.09%  "
Set-TextValue $ws.Range("D37") "0.0478"
Set-TextValue $ws.Range("E37") "  -1.54%  "
Set-TextValue $ws.Range("D38") "52.62"
Set-TextValue $ws.Range("E38") "  +2.46%  "
Set-TextValue $ws.Range("E39") "  +1.42%  "
Set-TextValue $ws.Range("E40") "  -0.21%  "
Set-TextValue $ws.Range("D41") "3.03"
Set-TextValue $ws.Range("E41") "  +0.69%  "
Set-TextValue $ws.Range("E42") "  +6.01%  "
Set-TextValue $ws.Range("D43") "0.125"
Set-TextValue $ws.Range("E43") "  +3.05%  "
Set-TextValue $ws.Range("D44") "134.55"
Set-TextValue $ws.Range("E44") "  -1.86%  "
Set-TextValue $ws.Range("D45") "17.04"
Set-TextValue $ws.Range("E45") "  +1.73%  "
Set-TextValue $ws.Range("D46") "0.282"
Set-TextValue $ws.Range("E46") "  -0.80%  "
Set-TextValue $ws.Range("D47") "3.88"
Set-TextValue $ws.Range("E47") "  -1.17%  "
Set-TextValue $ws.Range("D48") "2.19"
Set-TextValue $ws.Range("E48") "  -0.93%  "
Set-TextValue $ws.Range("D49") "21.89"
Set-TextValue $ws.Range("E49") "  -2.21%  "
Set-TextValue $ws.Range("D50") "2.205.25"
Set-TextValue $ws.Range("E50") "  +2.41%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D51") "3.768.77"
Set-TextValue $ws.Range("E51") "  +5.15%  "
